$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 18.74414781856099
$ws.Cells.Item(2, 3).Value = 11.41650246686654
$ws.Cells.Item(2, 5).Value = 15.77264841220765
$ws.Cells.Item(2, 6).Value = 37.67989913236096
$ws.Cells.Item(2, 7).Value = 3.651159446677739
$ws.Cells.Item(2, 10).Value = 7.897145562944408
$ws.Cells.Item(2, 12).Value = 12.74054999606443
$ws.Cells.Item(2, 15).Value = 24.11707720413692
$ws.Cells.Item(3, 2).Value = 18.15060801836043
$ws.Cells.Item(3, 3).Value = 11.25089673109048
$ws.Cells.Item(3, 5).Value = 15.80305615782843
$ws.Cells.Item(3, 6).Value = 37.76966863127084
$ws.Cells.Item(3, 7).Value = 3.653441036034238
$ws.Cells.Item(3, 10).Value = 7.905711441890567
$ws.Cells.Item(3, 12).Value = 12.69214155118813
$ws.Cells.Item(3, 15).Value = 24.24584427579515
$ws.Cells.Item(4, 2).Value = 17.77753823037522
$ws.Cells.Item(4, 3).Value = 11.14789228683985
$ws.Cells.Item(4, 5).Value = 15.82383742462421
$ws.Cells.Item(4, 6).Value = 37.83545521144573
$ws.Cells.Item(4, 7).Value = 3.654915341393064
$ws.Cells.Item(4, 10).Value = 7.911455398534524
$ws.Cells.Item(4, 12).Value = 12.66375553266892
$ws.Cells.Item(4, 15).Value = 24.33187061669377
$ws.Cells.Item(5, 2).Value = 17.6235546178762
$ws.Cells.Item(5, 3).Value = 11.10562049341983
$ws.Cells.Item(5, 5).Value = 15.83283701249318
$ws.Cells.Item(5, 6).Value = 37.86493712084933
$ws.Cells.Item(5, 7).Value = 3.655534650084171
$ws.Cells.Item(5, 10).Value = 7.91391828278268
$ws.Cells.Item(5, 12).Value = 12.65253120016222
$ws.Cells.Item(5, 15).Value = 24.36867131735724
$ws.Cells.Item(6, 2).Value = 17.59787437163662
$ws.Cells.Item(6, 3).Value = 11.09858443631823
$ws.Cells.Item(6, 5).Value = 15.83436346903716
$ws.Cells.Item(6, 6).Value = 37.86999374857911
$ws.Cells.Item(6, 7).Value = 3.655638605989306
$ws.Cells.Item(6, 10).Value = 7.914334632222275
$ws.Cells.Item(6, 12).Value = 12.65068833516502
$ws.Cells.Item(6, 15).Value = 24.37488720243772
$ws.Cells.Item(7, 2).Value = 17.77546916886692
$ws.Cells.Item(7, 3).Value = 11.14732334800171
$ws.Cells.Item(7, 5).Value = 15.82395664565134
$ws.Cells.Item(7, 6).Value = 37.83584200328936
$ws.Cells.Item(7, 7).Value = 3.654923618551716
$ws.Cells.Item(7, 10).Value = 7.911488118711802
$ws.Cells.Item(7, 12).Value = 12.66360275876897
$ws.Cells.Item(7, 15).Value = 24.33235987009516
$ws.Cells.Item(8, 2).Value = 18.54141163798687
$ws.Cells.Item(8, 3).Value = 11.3596971631183
$ws.Cells.Item(8, 5).Value = 15.78269508447835
$ws.Cells.Item(8, 6).Value = 37.70863232089724
$ws.Cells.Item(8, 7).Value = 3.651930940769566
$ws.Cells.Item(8, 10).Value = 7.899998735965022
$ws.Cells.Item(8, 12).Value = 12.72358557666103
$ws.Cells.Item(8, 15).Value = 24.16002699245087
$ws.Cells.Item(9, 2).Value = 19.96630938864775
$ws.Cells.Item(9, 3).Value = 11.76414681801424
$ws.Cells.Item(9, 5).Value = 15.71851666701057
$ws.Cells.Item(9, 6).Value = 37.5441900992102
$ws.Cells.Item(9, 7).Value = 3.64664195591925
$ws.Cells.Item(9, 10).Value = 7.881296803184569
$ws.Cells.Item(9, 12).Value = 12.85149954249576
$ws.Cells.Item(9, 15).Value = 23.87762403061025
$ws.Cells.Item(10, 2).Value = 20.95589401522366
$ws.Cells.Item(10, 3).Value = 12.05199903024433
$ws.Cells.Item(10, 5).Value = 15.68154904208
$ws.Cells.Item(10, 6).Value = 37.4756541834253
$ws.Cells.Item(10, 7).Value = 3.643105650006215
$ws.Cells.Item(10, 10).Value = 7.869870000069638
$ws.Cells.Item(10, 12).Value = 12.95130350110248
$ws.Cells.Item(10, 15).Value = 23.70440383411236
$ws.Cells.Item(11, 2).Value = 21.39181269595468
$ws.Cells.Item(11, 3).Value = 12.18053654162617
$ws.Cells.Item(11, 5).Value = 15.66693844230284
$ws.Cells.Item(11, 6).Value = 37.45590347369708
$ws.Cells.Item(11, 7).Value = 3.641571964239032
$ws.Cells.Item(11, 10).Value = 7.865169664990964
$ws.Cells.Item(11, 12).Value = 12.99787118792185
$ws.Cells.Item(11, 15).Value = 23.63312703189342
$ws.Cells.Item(12, 2).Value = 21.5547039257342
$ws.Cells.Item(12, 3).Value = 12.22883272593341
$ws.Cells.Item(12, 5).Value = 15.66172259399167
$ws.Cells.Item(12, 6).Value = 37.45007194035607
$ws.Cells.Item(12, 7).Value = 3.641001919557322
$ws.Cells.Item(12, 10).Value = 7.863460997389337
$ws.Cells.Item(12, 12).Value = 13.01566393979665
$ws.Cells.Item(12, 15).Value = 23.60722515342068
$ws.Cells.Item(13, 2).Value = 21.51972137364925
$ws.Cells.Item(13, 3).Value = 12.21844859134715
$ws.Cells.Item(13, 5).Value = 15.66283183383715
$ws.Cells.Item(13, 6).Value = 37.45125452714377
$ws.Cells.Item(13, 7).Value = 3.641124212594504
$ws.Cells.Item(13, 10).Value = 7.86382582606514
$ws.Cells.Item(13, 12).Value = 13.0118250472961
$ws.Cells.Item(13, 15).Value = 23.61275503009545
$ws.Cells.Item(14, 2).Value = 21.40525828796305
$ws.Cells.Item(14, 3).Value = 12.18451765488149
$ws.Cells.Item(14, 5).Value = 15.66650298309495
$ws.Cells.Item(14, 6).Value = 37.45539066726605
$ws.Cells.Item(14, 7).Value = 3.641524851616385
$ws.Cells.Item(14, 10).Value = 7.865027665877109
$ws.Cells.Item(14, 12).Value = 12.99933188328208
$ws.Cells.Item(14, 15).Value = 23.6309741990219
$ws.Cells.Item(15, 2).Value = 21.33485853301488
$ws.Cells.Item(15, 3).Value = 12.16368380300221
$ws.Cells.Item(15, 5).Value = 15.66879292185629
$ws.Cells.Item(15, 6).Value = 37.45813885671387
$ws.Cells.Item(15, 7).Value = 3.641771650148566
$ws.Cells.Item(15, 10).Value = 7.865773096119952
$ws.Cells.Item(15, 12).Value = 12.99169983632727
$ws.Cells.Item(15, 15).Value = 23.64227601957253
$ws.Cells.Item(16, 2).Value = 20.92710729266929
$ws.Cells.Item(16, 3).Value = 12.04354758664582
$ws.Cells.Item(16, 5).Value = 15.68254824165971
$ws.Cells.Item(16, 6).Value = 37.47717528214687
$ws.Cells.Item(16, 7).Value = 3.64320738434733
$ws.Cells.Item(16, 10).Value = 7.870187166112054
$ws.Cells.Item(16, 12).Value = 12.94828288698921
$ws.Cells.Item(16, 15).Value = 23.70921393937345
$ws.Cells.Item(17, 2).Value = 20.67321644375579
$ws.Cells.Item(17, 3).Value = 11.96920802157477
$ws.Cells.Item(17, 5).Value = 15.6915514819906
$ws.Cells.Item(17, 6).Value = 37.49178370616863
$ws.Cells.Item(17, 7).Value = 3.644107329671389
$ws.Cells.Item(17, 10).Value = 7.873022310582462
$ws.Cells.Item(17, 12).Value = 12.92194029473816
$ws.Cells.Item(17, 15).Value = 23.75221004794584
$ws.Cells.Item(18, 2).Value = 20.52585096677067
$ws.Cells.Item(18, 3).Value = 11.92622504831341
$ws.Cells.Item(18, 5).Value = 15.6969375867736
$ws.Cells.Item(18, 6).Value = 37.5012614153811
$ws.Cells.Item(18, 7).Value = 3.644632016846905
$ws.Cells.Item(18, 10).Value = 7.874699883405048
$ws.Cells.Item(18, 12).Value = 12.906899078045
$ws.Cells.Item(18, 15).Value = 23.77764782838442
$ws.Cells.Item(19, 2).Value = 20.47573069777866
$ws.Cells.Item(19, 3).Value = 11.91163413484224
$ws.Cells.Item(19, 5).Value = 15.69879690909841
$ws.Cells.Item(19, 6).Value = 37.5046549381909
$ws.Cells.Item(19, 7).Value = 3.644810881592425
$ws.Cells.Item(19, 10).Value = 7.875275941315054
$ws.Cells.Item(19, 12).Value = 12.90182560742978
$ws.Cells.Item(19, 15).Value = 23.7863819195234
$ws.Cells.Item(20, 2).Value = 20.70038265010196
$ws.Cells.Item(20, 3).Value = 11.97714509918277
$ws.Cells.Item(20, 5).Value = 15.69057158026257
$ws.Cells.Item(20, 6).Value = 37.49011728561568
$ws.Cells.Item(20, 7).Value = 3.644010798425512
$ws.Cells.Item(20, 10).Value = 7.87271565586586
$ws.Cells.Item(20, 12).Value = 12.92473315383375
$ws.Cells.Item(20, 15).Value = 23.74755975523854
$ws.Cells.Item(21, 2).Value = 21.43893906641921
$ws.Cells.Item(21, 3).Value = 12.19449450449952
$ws.Cells.Item(21, 5).Value = 15.66541608103335
$ws.Cells.Item(21, 6).Value = 37.45413103624581
$ws.Cells.Item(21, 7).Value = 3.641406883488789
$ws.Cells.Item(21, 10).Value = 7.864672725402776
$ws.Cells.Item(21, 12).Value = 13.00299719850863
$ws.Cells.Item(21, 15).Value = 23.62559317039525
$ws.Cells.Item(22, 2).Value = 21.90886039227076
$ws.Cells.Item(22, 3).Value = 12.33432678672632
$ws.Cells.Item(22, 5).Value = 15.65082225788868
$ws.Cells.Item(22, 6).Value = 37.44021687048983
$ws.Cells.Item(22, 7).Value = 3.639767584814188
$ws.Cells.Item(22, 10).Value = 7.859831335582783
$ws.Cells.Item(22, 12).Value = 13.05506721289947
$ws.Cells.Item(22, 15).Value = 23.55223249203947
$ws.Cells.Item(23, 2).Value = 21.65926287657541
$ws.Cells.Item(23, 3).Value = 12.25990879891345
$ws.Cells.Item(23, 5).Value = 15.65844241280892
$ws.Cells.Item(23, 6).Value = 37.44676308332716
$ws.Cells.Item(23, 7).Value = 3.640636807985318
$ws.Cells.Item(23, 10).Value = 7.862377401530122
$ws.Cells.Item(23, 12).Value = 13.02719536620524
$ws.Cells.Item(23, 15).Value = 23.59080284834499
$ws.Cells.Item(24, 2).Value = 20.68810515942497
$ws.Cells.Item(24, 3).Value = 11.97355750451018
$ws.Cells.Item(24, 5).Value = 15.69101393961195
$ws.Cells.Item(24, 6).Value = 37.49086731323476
$ws.Cells.Item(24, 7).Value = 3.644054417468188
$ws.Cells.Item(24, 10).Value = 7.872854146121626
$ws.Cells.Item(24, 12).Value = 12.92347017919441
$ws.Cells.Item(24, 15).Value = 23.74965991422269
$ws.Cells.Item(25, 2).Value = 19.59022105817964
$ws.Cells.Item(25, 3).Value = 11.6562266597933
$ws.Cells.Item(25, 5).Value = 15.73408871360772
$ws.Cells.Item(25, 6).Value = 37.57952395181037
$ws.Cells.Item(25, 7).Value = 3.648011111400145
$ws.Cells.Item(25, 10).Value = 7.885948412413112
$ws.Cells.Item(25, 12).Value = 12.81583893985806
$ws.Cells.Item(25, 15).Value = 23.94803391954737
